$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Runmode column (C3:C7) from "N" to "Y"
$ws.Range("C3:C7").Value = "Y"

# Make this sheet active and select the updated range, matching the new selection
$ws.Activate()
$ws.Range("C2:C7").Select()
